# mise à jour nb de mux necessaire
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Update the quantity of the "Demux 4 vers 1" component (row 6) from 4 to 6
$ws.Range("C6").Value = 6

# Move the active selection (matches the saved cursor position in the diff)
$ws.Range("I10").Select()
